$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (odd/primary) Pearson logo: id=1, image2.png -> image1.png ---
$f1 = $sec.Footers.Item(1)
$shp1 = $f1.Range.InlineShapes.Item(1).ConvertToShape()
$shp1.Name = "image1.png"
$shp1.ConvertToInlineShape() | Out-Null

# --- Footer (even) Pearson logo: id=2, image2.png -> image1.png ---
$f2 = $sec.Footers.Item(2)
$shp2 = $f2.Range.InlineShapes.Item(1).ConvertToShape()
$shp2.Name = "image1.png"
$shp2.ConvertToInlineShape() | Out-Null

# --- Header (even) BTec logo: id=3, image1.jpg -> image2.jpg ---
$h2 = $sec.Headers.Item(2)
$shp3 = $h2.Range.InlineShapes.Item(1).ConvertToShape()
$shp3.Name = "image2.jpg"
$shp3.ConvertToInlineShape() | Out-Null
